# Applies the 2025-07-14 09:30:09 automatic update to the workbook.
$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
# Row 19 (LOZANO MOLINA TITO / RENOVA&DISEÑA S.A.)
$wsVentasPorGrupo.Range("L19").Value = 309.47
$wsVentasPorGrupo.Range("M19").Value = 1735.84

# Row 29 (totals row, counters formatted as "N de 27")
$wsVentasPorGrupo.Range("L29").Value = "1 de 27"
$wsVentasPorGrupo.Range("M29").Value = "2 de 27"

# --- Sheet "VENTA MENSUAL" ---
# Row 19 (LOZANO MOLINA TITO / RENOVA&DISEÑA S.A.)
$wsVentaMensual.Range("F19").Value = 2045.31

# Row 29 (totals row)
$wsVentaMensual.Range("F29").Value = 3664.27

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
# Row 15 (PIEDRA SINTERIZADA)
$wsCumplimiento.Range("D15").Value = 309.47
$wsCumplimiento.Range("E15").Value = 1328.53
$wsCumplimiento.Range("F15").Value = 0.1889316239316239

# Row 16 (PORCELANATO)
$wsCumplimiento.Range("D16").Value = 3913.12
$wsCumplimiento.Range("E16").Value = 13172.77
$wsCumplimiento.Range("F16").Value = 0.2290264071698928

# Row 19 (TOTAL)
$wsCumplimiento.Range("D19").Value = 3664.27
$wsCumplimiento.Range("E19").Value = 23517.04093005039
$wsCumplimiento.Range("F19").Value = 0.1348084354514688
